$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.796899999999999
$ws.Range("A3").Value = -21.38920000000002
$ws.Range("B5").Value = 4.882200000000002
$ws.Range("C5").Value = -14.0346
$ws.Range("D7").Value = -8.012599999999997
$ws.Range("E7").Value = 12.953
$ws.Range("C9").Value = -11.85300000000002
$ws.Range("C11").Value = -13.1283
$ws.Range("D11").Value = -8.044499999999998
$ws.Range("A14").Value = -20.52079999999998
$ws.Range("D19").Value = -8.257699999999996
$ws.Range("A21").Value = -21.17290000000001
$ws.Range("C21").Value = -10.59689999999999
$ws.Range("D21").Value = -7.671499999999996
$ws.Range("A23").Value = -21.30490000000002
$ws.Range("A25").Value = -22.51290000000003
